$d = $word.ActiveDocument
$sec = $d.Sections(1)
$f = $sec.Footers(1)
$tbl = $d.Tables.Add($f.Range, 1, 3)
$tblRange = $tbl.Range
$tblRange.Collapse(0)
$newPara = $f.Range.Paragraphs.Add($tblRange)
Write-Output "count after add: $($f.Range.Paragraphs.Count)"
